# Update the "想去人数" (F column) counts that changed in the source data
# refresh across all four sheets: 展览, 演出, 本地生活, 全部类型.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1745
$ws.Range("F5").Value = 439
$ws.Range("F6").Value = 813
$ws.Range("F7").Value = 243
$ws.Range("F8").Value = 1169
$ws.Range("F9").Value = 333
$ws.Range("F11").Value = 873
$ws.Range("F12").Value = 677
$ws.Range("F13").Value = 183
$ws.Range("F14").Value = 504
$ws.Range("F15").Value = 140
$ws.Range("F17").Value = 168
$ws.Range("F18").Value = 2903
$ws.Range("F19").Value = 2614
$ws.Range("F26").Value = 5255
$ws.Range("F27").Value = 589
$ws.Range("F28").Value = 976
$ws.Range("F29").Value = 21
$ws.Range("F31").Value = 299
$ws.Range("F32").Value = 1084
$ws.Range("F33").Value = 70
$ws.Range("F35").Value = 286
$ws.Range("F36").Value = 36

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1113
$ws.Range("F10").Value = 28
$ws.Range("F14").Value = 607
$ws.Range("F17").Value = 983
$ws.Range("F19").Value = 41
$ws.Range("F24").Value = 310
$ws.Range("F25").Value = 276
$ws.Range("F26").Value = 3908
$ws.Range("F29").Value = 20
$ws.Range("F31").Value = 51

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1779
$ws.Range("F5").Value = 2444
$ws.Range("F6").Value = 1029
$ws.Range("F7").Value = 3
$ws.Range("F9").Value = 1309
$ws.Range("F10").Value = 357

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1779
$ws.Range("F4").Value = 2444
$ws.Range("F5").Value = 1745
$ws.Range("F6").Value = 1029
$ws.Range("F7").Value = 1309
$ws.Range("F8").Value = 357
$ws.Range("F11").Value = 439
$ws.Range("F12").Value = 813
$ws.Range("F13").Value = 243
$ws.Range("F14").Value = 1169
$ws.Range("F15").Value = 333
$ws.Range("F16").Value = 873
$ws.Range("F17").Value = 677
$ws.Range("F18").Value = 1113
$ws.Range("F19").Value = 1113
$ws.Range("F20").Value = 504
$ws.Range("F22").Value = 168
$ws.Range("F23").Value = 2903
$ws.Range("F24").Value = 2614
$ws.Range("F28").Value = 28
$ws.Range("F30").Value = 5255
$ws.Range("F31").Value = 589
$ws.Range("F32").Value = 976
$ws.Range("F33").Value = 607
$ws.Range("F34").Value = 607
$ws.Range("F35").Value = 21
$ws.Range("F38").Value = 299
$ws.Range("F41").Value = 41
$ws.Range("F44").Value = 310
$ws.Range("F45").Value = 310
$ws.Range("F46").Value = 276
$ws.Range("F47").Value = 1084
$ws.Range("F51").Value = 286
$ws.Range("F52").Value = 36
